$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2024-07-30 Tuesday"

# Update the 20x5 grid of arithmetic expressions in the table
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "30+10="
$t.Cell(1,2).Range.Text = "37-22="
$t.Cell(1,3).Range.Text = "88-48="
$t.Cell(1,4).Range.Text = "50-22="
$t.Cell(1,5).Range.Text = "63-38="

$t.Cell(2,1).Range.Text = "24-22="
$t.Cell(2,2).Range.Text = "11-11="
$t.Cell(2,3).Range.Text = "66-22="
$t.Cell(2,4).Range.Text = "19-6="
$t.Cell(2,5).Range.Text = "48-23="

$t.Cell(3,1).Range.Text = "99-88="
$t.Cell(3,2).Range.Text = "27-26="
$t.Cell(3,3).Range.Text = "39-35="
$t.Cell(3,4).Range.Text = "81-8="
$t.Cell(3,5).Range.Text = "5+43="

$t.Cell(4,1).Range.Text = "83-29="
$t.Cell(4,2).Range.Text = "34+36="
$t.Cell(4,3).Range.Text = "91-32="
$t.Cell(4,4).Range.Text = "27+37="
$t.Cell(4,5).Range.Text = "7+59="

$t.Cell(5,1).Range.Text = "59+0="
$t.Cell(5,2).Range.Text = "46+21="
$t.Cell(5,3).Range.Text = "12+18="
$t.Cell(5,4).Range.Text = "87+10="
$t.Cell(5,5).Range.Text = "20+28="

$t.Cell(6,1).Range.Text = "25+60="
$t.Cell(6,2).Range.Text = "24+66="
$t.Cell(6,3).Range.Text = "31+27="
$t.Cell(6,4).Range.Text = "32-25="
$t.Cell(6,5).Range.Text = "56+4="

$t.Cell(7,1).Range.Text = "47-18="
$t.Cell(7,2).Range.Text = "42-27="
$t.Cell(7,3).Range.Text = "94-76="
$t.Cell(7,4).Range.Text = "34+47="
$t.Cell(7,5).Range.Text = "71-36="

$t.Cell(8,1).Range.Text = "32-21="
$t.Cell(8,2).Range.Text = "13+67="
$t.Cell(8,3).Range.Text = "72+3="
$t.Cell(8,4).Range.Text = "81-49="
$t.Cell(8,5).Range.Text = "77-20="

$t.Cell(9,1).Range.Text = "23+67="
$t.Cell(9,2).Range.Text = "83-54="
$t.Cell(9,3).Range.Text = "71-52="
$t.Cell(9,4).Range.Text = "57-55="
$t.Cell(9,5).Range.Text = "95-60="

$t.Cell(10,1).Range.Text = "8+45="
$t.Cell(10,2).Range.Text = "36+30="
$t.Cell(10,3).Range.Text = "60-12="
$t.Cell(10,4).Range.Text = "91+6="
$t.Cell(10,5).Range.Text = "22+49="

$t.Cell(11,1).Range.Text = "27-23="
$t.Cell(11,2).Range.Text = "5+50="
$t.Cell(11,3).Range.Text = "6-0="
$t.Cell(11,4).Range.Text = "12-10="
$t.Cell(11,5).Range.Text = "30+8="

$t.Cell(12,1).Range.Text = "23+0="
$t.Cell(12,2).Range.Text = "29+29="
$t.Cell(12,3).Range.Text = "96-77="
$t.Cell(12,4).Range.Text = "78-40="
$t.Cell(12,5).Range.Text = "32+2="

$t.Cell(13,1).Range.Text = "60-31="
$t.Cell(13,2).Range.Text = "33+38="
$t.Cell(13,3).Range.Text = "40+44="
$t.Cell(13,4).Range.Text = "3+22="
$t.Cell(13,5).Range.Text = "65-26="

$t.Cell(14,1).Range.Text = "75+10="
$t.Cell(14,2).Range.Text = "45+47="
$t.Cell(14,3).Range.Text = "79-77="
$t.Cell(14,4).Range.Text = "45+8="
$t.Cell(14,5).Range.Text = "49+0="

$t.Cell(15,1).Range.Text = "60-0="
$t.Cell(15,2).Range.Text = "59+36="
$t.Cell(15,3).Range.Text = "64-38="
$t.Cell(15,4).Range.Text = "45-0="
$t.Cell(15,5).Range.Text = "40-14="

$t.Cell(16,1).Range.Text = "28-3="
$t.Cell(16,2).Range.Text = "6+88="
$t.Cell(16,3).Range.Text = "2+74="
$t.Cell(16,4).Range.Text = "77-43="
$t.Cell(16,5).Range.Text = "8+12="

$t.Cell(17,1).Range.Text = "53+9="
$t.Cell(17,2).Range.Text = "20+30="
$t.Cell(17,3).Range.Text = "40+56="
$t.Cell(17,4).Range.Text = "54+29="
$t.Cell(17,5).Range.Text = "92-77="

$t.Cell(18,1).Range.Text = "79-60="
$t.Cell(18,2).Range.Text = "54-33="
$t.Cell(18,3).Range.Text = "33-21="
$t.Cell(18,4).Range.Text = "46+42="
$t.Cell(18,5).Range.Text = "87-50="

$t.Cell(19,1).Range.Text = "59-0="
$t.Cell(19,2).Range.Text = "2+87="
$t.Cell(19,3).Range.Text = "65-49="
$t.Cell(19,4).Range.Text = "99-10="
$t.Cell(19,5).Range.Text = "11+82="

$t.Cell(20,1).Range.Text = "74-9="
$t.Cell(20,2).Range.Text = "31+40="
$t.Cell(20,3).Range.Text = "47-18="
$t.Cell(20,4).Range.Text = "0+57="
$t.Cell(20,5).Range.Text = "57+30="
